$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before the current Policy_Id column (C) so that
# Policy_Id / Policy_Status shift right to D / E, making room for the
# new "Employee_Code" column.
$ws.Columns("C:C").Insert()

# New header for the inserted column.
$ws.Range("C1").Value = "Employee_Code"

# Remove the now-obsolete trailing data rows (old rows 5, 6, 7).
$ws.Rows("5:7").Delete()

# Replace the remaining data rows with the updated records.
$ws.Range("A2").Value = "Jayashree Kulai"
$ws.Range("B2").Value = "jayashree.cs16@sahyadri.edu.in"
$ws.Range("C2").Value = "MNG001"
$ws.Range("D2").Value = "5fd0709530a434204c3007d5"
$ws.Range("E2").Value = $false

$ws.Range("A3").Value = "Roy Pashan"
$ws.Range("B3").Value = "mail2winstonroy@yahoo.com"
$ws.Range("C3").Value = "MNG002"
$ws.Range("D3").Value = "5fd0709530a434204c3007d5"
$ws.Range("E3").Value = $false

$ws.Range("A4").Value = "Winston Roy"
$ws.Range("B4").Value = "pashanwinsty1998@gmail.com"
$ws.Range("C4").Value = "AD002"
$ws.Range("D4").Value = "5fd0709530a434204c3007d5"
$ws.Range("E4").Value = $false
